$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from G1 (bold/border/centered) onto new header cell H1
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# New "Save" data column values
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 1
